$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force affected range to Text format first so numeric-looking values
# (e.g. "1.003", "313.66") are preserved as text, matching the source data shape.
$ws.Range("B2:E51").NumberFormat = "@"

$ws.Range('D2').Value = '28.089.35'
$ws.Range('E2').Value = '  -0.41%  '
$ws.Range('D3').Value = '1.876.73'
$ws.Range('E3').Value = '  -1.78%  '
$ws.Range('D4').Value = '1.003'
$ws.Range('E4').Value = '  +0.17%  '
$ws.Range('D5').Value = '313.66'
$ws.Range('E6').Value = '  +0.17%  '
$ws.Range('E7').Value = '  -0.34%  '
$ws.Range('D8').Value = '0.3839'
$ws.Range('E8').Value = '  -2.02%  '
$ws.Range('D9').Value = '0.08638'
$ws.Range('E9').Value = '  -7.03%  '
$ws.Range('E10').Value = '  -1.98%  '
$ws.Range('D11').Value = '41.48'
$ws.Range('E11').Value = '  -0.95%  '
$ws.Range('D12').Value = '6.333'
$ws.Range('E12').Value = '  -0.86%  '
$ws.Range('D13').Value = '20.68'
$ws.Range('E13').Value = '  -1.15%  '
$ws.Range('D14').Value = '1.880.76'
$ws.Range('E14').Value = '  -1.63%  '
$ws.Range('D15').Value = '1.003'
$ws.Range('E15').Value = '  +0.13%  '
$ws.Range('D16').Value = '7.174'
$ws.Range('E16').Value = '  -1.95%  '
$ws.Range('D17').Value = '0.00001103'
$ws.Range('E17').Value = '  -1.82%  '
$ws.Range('D18').Value = '91.10'
$ws.Range('E18').Value = '  -1.47%  '
$ws.Range('D19').Value = '0.06637'
$ws.Range('D20').Value = '18.14'
$ws.Range('E20').Value = '  +0.86%  '
$ws.Range('E21').Value = '  +0.23%  '
$ws.Range('D22').Value = '6.105'
$ws.Range('E22').Value = '  -1.88%  '
$ws.Range('D23').Value = '28.129.18'
$ws.Range('E23').Value = '  -0.51%  '
$ws.Range('D24').Value = '11.45'
$ws.Range('E24').Value = '  +0.23%  '
$ws.Range('D25').Value = '2.267'
$ws.Range('E25').Value = '  -2.24%  '
$ws.Range('B26').Value = 'LEO'
$ws.Range('C26').Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range('D26').Value = '3.402'
$ws.Range('E26').Value = '  +0.25%  '
$ws.Range('B27').Value = 'LidoDAOToken'
$ws.Range('C27').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D27').Value = '2.569'
$ws.Range('E27').Value = '  -0.87%  '
$ws.Range('B28').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C28').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D28').Value = '2.084.97'
$ws.Range('E28').Value = '  -1.83%  '
$ws.Range('B29').Value = 'EthereumClassic'
$ws.Range('C29').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D29').Value = '20.73'
$ws.Range('E29').Value = '  -1.59%  '
$ws.Range('B30').Value = 'Monero'
$ws.Range('C30').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D30').Value = '157.20'
$ws.Range('E30').Value = '  -0.43%  '
$ws.Range('B31').Value = 'BitcoinCash'
$ws.Range('C31').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D31').Value = '125.97'
$ws.Range('E31').Value = '  -0.99%  '
$ws.Range('B32').Value = 'Stellar'
$ws.Range('C32').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D32').Value = '0.1050'
$ws.Range('E32').Value = '  -2.30%  '
$ws.Range('B33').Value = 'ImmutableX'
$ws.Range('C33').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D33').Value = '1.061'
$ws.Range('E33').Value = '  -3.44%  '
$ws.Range('B34').Value = 'Filecoin'
$ws.Range('C34').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D34').Value = '5.595'
$ws.Range('E34').Value = '  -0.72%  '
$ws.Range('B35').Value = 'HuobiToken'
$ws.Range('C35').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D35').Value = '3.602'
$ws.Range('E35').Value = '  -0.39%  '
$ws.Range('B36').Value = 'FraxShare'
$ws.Range('C36').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D36').Value = '9.700'
$ws.Range('E36').Value = '  +0.52%  '
$ws.Range('B37').Value = 'VeChain'
$ws.Range('C37').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D37').Value = '0.02456'
$ws.Range('E37').Value = '  +1.40%  '
$ws.Range('B38').Value = 'Hedera'
$ws.Range('C38').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D38').Value = '0.06579'
$ws.Range('E38').Value = '  -1.27%  '
$ws.Range('B39').Value = 'Algorand'
$ws.Range('C39').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D39').Value = '0.2174'
$ws.Range('E39').Value = '  -1.10%  '
$ws.Range('B40').Value = 'ARBITRUM'
$ws.Range('C40').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D40').Value = '1.205'
$ws.Range('E40').Value = '  -3.32%  '
$ws.Range('B41').Value = 'TrustWalletToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D41').Value = '1.246'
$ws.Range('E41').Value = '  -2.71%  '
$ws.Range('B42').Value = 'Aptos'
$ws.Range('C42').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D42').Value = '11.58'
$ws.Range('E42').Value = '  +0.83%  '
$ws.Range('B43').Value = 'TheSandbox'
$ws.Range('C43').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D43').Value = '0.6374'
$ws.Range('E43').Value = '  -1.11%  '
$ws.Range('B44').Value = 'InternetComputer(DFINITY)'
$ws.Range('C44').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D44').Value = '4.899'
$ws.Range('E44').Value = '  -1.98%  '
$ws.Range('B45').Value = 'EnergySwap'
$ws.Range('C45').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D45').Value = '13.15'
$ws.Range('E45').Value = '  -0.93%  '
$ws.Range('B46').Value = 'Decentraland'
$ws.Range('C46').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D46').Value = '0.5987'
$ws.Range('E46').Value = '  -0.74%  '
$ws.Range('B47').Value = 'WEMIXTOKEN'
$ws.Range('C47').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D47').Value = '1.280'
$ws.Range('E47').Value = '  -0.12%  '
$ws.Range('B48').Value = 'PancakeSwap'
$ws.Range('C48').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D48').Value = '3.673'
$ws.Range('E48').Value = '  -1.31%  '
$ws.Range('B49').Value = 'EOS'
$ws.Range('C49').Value = 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos'
$ws.Range('D49').Value = '1.230'
$ws.Range('E49').Value = '  +3.80%  '
$ws.Range('B50').Value = 'NEARProtocol'
$ws.Range('C50').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D50').Value = '1.991'
$ws.Range('E50').Value = '  -1.59%  '
$ws.Range('B51').Value = 'Quant'
$ws.Range('C51').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D51').Value = '121.31'
$ws.Range('E51').Value = '  -1.34%  '

# Restore default (Normal) style so no extra formatting/number-format is visibly applied
$ws.Range("B2:E51").Style = "Normal"

